$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I12").Value = "sd"
$ws.Range("J12").Value = "Statement-non-opinion"
$ws.Range("I33").Value = "aa"
$ws.Range("J33").Value = "Agree/Accept"
$ws.Range("I37").Value = "sd"
$ws.Range("J37").Value = "Statement-non-opinion"
$ws.Range("I44").Value = "%"
$ws.Range("J44").Value = "Uninterpretable"
$ws.Range("I82").Value = "sd"
$ws.Range("J82").Value = "Statement-non-opinion"
$ws.Range("I84").Value = "sv"
$ws.Range("J84").Value = "Statement-opinion"
$ws.Range("I91").Value = "sd"
$ws.Range("J91").Value = "Statement-non-opinion"
$ws.Range("I101").Value = "sd"
$ws.Range("J101").Value = "Statement-non-opinion"
$ws.Range("I104").Value = "aa"
$ws.Range("J104").Value = "Agree/Accept"
$ws.Range("I114").Value = "sv"
$ws.Range("J114").Value = "Statement-opinion"
$ws.Range("I128").Value = "sd"
$ws.Range("J128").Value = "Statement-non-opinion"
$ws.Range("I186").Value = "sv"
$ws.Range("J186").Value = "Statement-opinion"
$ws.Range("I194").Value = "aa"
$ws.Range("J194").Value = "Agree/Accept"
$ws.Range("I203").Value = "aa"
$ws.Range("J203").Value = "Agree/Accept"
$ws.Range("I214").Value = "aa"
$ws.Range("J214").Value = "Agree/Accept"
$ws.Range("I217").Value = "sd"
$ws.Range("J217").Value = "Statement-non-opinion"
$ws.Range("I229").Value = "%"
$ws.Range("J229").Value = "Uninterpretable"
$ws.Range("I243").Value = "%"
$ws.Range("J243").Value = "Uninterpretable"
$ws.Range("I256").Value = "sv"
$ws.Range("J256").Value = "Statement-opinion"
$ws.Range("I263").Value = "ba"
$ws.Range("J263").Value = "Appreciation"
$ws.Range("I266").Value = "sd"
$ws.Range("J266").Value = "Statement-non-opinion"
$ws.Range("I289").Value = "%"
$ws.Range("J289").Value = "Uninterpretable"
$ws.Range("I290").Value = "aa"
$ws.Range("J290").Value = "Agree/Accept"
$ws.Range("I322").Value = "sd"
$ws.Range("J322").Value = "Statement-non-opinion"
$ws.Range("I323").Value = "sd"
$ws.Range("J323").Value = "Statement-non-opinion"
$ws.Range("I326").Value = "sd"
$ws.Range("J326").Value = "Statement-non-opinion"
$ws.Range("I342").Value = "%"
$ws.Range("J342").Value = "Uninterpretable"
$ws.Range("I350").Value = "b"
$ws.Range("J350").Value = "Acknowledge (Backchannel)"
$ws.Range("I352").Value = "aa"
$ws.Range("J352").Value = "Agree/Accept"
$ws.Range("I364").Value = "ba"
$ws.Range("J364").Value = "Appreciation"
$ws.Range("I370").Value = "sv"
$ws.Range("J370").Value = "Statement-opinion"
$ws.Range("I376").Value = "sd"
$ws.Range("J376").Value = "Statement-non-opinion"
$ws.Range("I378").Value = "sd"
$ws.Range("J378").Value = "Statement-non-opinion"
$ws.Range("I380").Value = "%"
$ws.Range("J380").Value = "Uninterpretable"
$ws.Range("I389").Value = "sd"
$ws.Range("J389").Value = "Statement-non-opinion"
$ws.Range("I397").Value = "sd"
$ws.Range("J397").Value = "Statement-non-opinion"
$ws.Range("I400").Value = "ba"
$ws.Range("J400").Value = "Appreciation"
$ws.Range("I409").Value = "ba"
$ws.Range("J409").Value = "Appreciation"
$ws.Range("I413").Value = "sv"
$ws.Range("J413").Value = "Statement-opinion"
$ws.Range("I414").Value = "aa"
$ws.Range("J414").Value = "Agree/Accept"
$ws.Range("I420").Value = "sd"
$ws.Range("J420").Value = "Statement-non-opinion"
$ws.Range("I431").Value = "aa"
$ws.Range("J431").Value = "Agree/Accept"
$ws.Range("I432").Value = "sv"
$ws.Range("J432").Value = "Statement-opinion"
$ws.Range("I442").Value = "sd"
$ws.Range("J442").Value = "Statement-non-opinion"
$ws.Range("I444").Value = "sv"
$ws.Range("J444").Value = "Statement-opinion"
$ws.Range("I449").Value = "aa"
$ws.Range("J449").Value = "Agree/Accept"
$ws.Range("I456").Value = "aa"
$ws.Range("J456").Value = "Agree/Accept"
$ws.Range("I457").Value = "%"
$ws.Range("J457").Value = "Uninterpretable"
$ws.Range("I458").Value = "%"
$ws.Range("J458").Value = "Uninterpretable"
$ws.Range("I461").Value = "sd"
$ws.Range("J461").Value = "Statement-non-opinion"
$ws.Range("I464").Value = "aa"
$ws.Range("J464").Value = "Agree/Accept"
